# The deck ships two embedded DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" color scheme (used by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral" color scheme      (used by the Slide Master,
#                                                          i.e. the design applied to
#                                                          every slide)
# The authored edit swaps the two themes' contents: the deck's live/visible design
# (the Slide Master's theme) switches from the "Integral" palette to the "Office"
# palette. Font scheme and format scheme (fills/lines/effects) are identical between
# the two themes, so only the 12 theme colors change.
#
# Reproduce this by recoloring the Slide Master's theme color scheme, slot by slot,
# to the "Office" palette (the values that used to live in theme1.xml / "Office
# Theme").

function ConvertTo-BgrColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b -shl 16) -bor ($g -shl 8) -bor $r
}

# Office Theme color scheme, in PowerPoint's ColorScheme.Colors(1..12) order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6,
# hlink, folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $slot = $colorScheme.Colors($i)
    $slot.RGB = ConvertTo-BgrColor $officeThemeColors[$i - 1]
}

# Best-effort: try to restore the design/theme display name too (no-op on hosts
# that treat Design.Name as read-only, but harmless).
try {
    $design = $p.Designs.Item(1)
    $design.Name = "Office Theme"
} catch {
}
